# Update cryptocurrency price (col D) and 1h volume change (col E) values
# to match latest scrape, preserving the original text-literal cell storage
# (NumberFormat "@" forces text entry for numeric-looking strings like "1.001"
# or "21.00" so Excel does not silently coerce them to numbers; ClearFormats()
# afterwards drops the now-unneeded style so cell styling is left untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "29.724.50"
Set-TextValue "E2" "  +4.06%  "
Set-TextValue "D3" "1.913.43"
Set-TextValue "E3" "  +2.03%  "
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  -0.91%  "
Set-TextValue "D5" "318.63"
Set-TextValue "E5" "  +0.99%  "
Set-TextValue "E6" "  -0.78%  "
Set-TextValue "D7" "0.5197"
Set-TextValue "E7" "  +1.92%  "
Set-TextValue "D8" "0.3976"
Set-TextValue "E8" "  +1.30%  "
Set-TextValue "E9" "  +1.38%  "
Set-TextValue "E10" "  +2.75%  "
Set-TextValue "E11" "  +1.15%  "
Set-TextValue "D12" "6.300"
Set-TextValue "E12" "  +0.47%  "
Set-TextValue "D13" "1.908.06"
Set-TextValue "E13" "  +1.43%  "
Set-TextValue "D14" "20.94"
Set-TextValue "E14" "  +2.21%  "
Set-TextValue "D15" "7.359"
Set-TextValue "E15" "  +1.23%  "
Set-TextValue "D16" "1.002"
Set-TextValue "E16" "  -0.63%  "
Set-TextValue "D17" "94.06"
Set-TextValue "E17" "  +2.87%  "
Set-TextValue "D18" "0.00001116"
Set-TextValue "E18" "  +0.84%  "
Set-TextValue "D19" "0.06748"
Set-TextValue "E19" "  +0.14%  "
Set-TextValue "E20" "  +1.34%  "
Set-TextValue "E21" "  -0.81%  "
Set-TextValue "D22" "6.034"
Set-TextValue "E22" "  +1.21%  "
Set-TextValue "D23" "29.711.16"
Set-TextValue "E23" "  +3.91%  "
Set-TextValue "D25" "2.209"
Set-TextValue "E25" "  -1.65%  "
Set-TextValue "D26" "2.128.45"
Set-TextValue "E26" "  +1.53%  "
Set-TextValue "D27" "21.00"
Set-TextValue "E27" "  +1.17%  "
Set-TextValue "D28" "159.36"
Set-TextValue "E28" "  -1.57%  "
Set-TextValue "D29" "2.472"
Set-TextValue "E29" "  +4.48%  "
Set-TextValue "D30" "128.45"
Set-TextValue "E30" "  +1.43%  "
Set-TextValue "E31" "  +3.30%  "
Set-TextValue "D32" "0.1060"
Set-TextValue "E32" "  +0.40%  "
Set-TextValue "D33" "6.188"
Set-TextValue "E33" "  +6.60%  "
Set-TextValue "D34" "3.682"
Set-TextValue "E34" "  +1.99%  "
Set-TextValue "D35" "0.02502"
Set-TextValue "E35" "  +1.85%  "
Set-TextValue "D36" "0.06647"
Set-TextValue "E36" "  +1.73%  "
Set-TextValue "D37" "9.173"
Set-TextValue "E37" "  +3.04%  "
Set-TextValue "E38" "  +1.17%  "
Set-TextValue "D39" "1.245"
Set-TextValue "E39" "  +4.16%  "
Set-TextValue "D40" "5.199"
Set-TextValue "E40" "  +2.61%  "
Set-TextValue "D41" "0.6544"
Set-TextValue "E41" "  +1.21%  "
Set-TextValue "D42" "1.242"
Set-TextValue "E42" "  -1.71%  "
Set-TextValue "D43" "11.48"
Set-TextValue "E43" "  +2.59%  "
Set-TextValue "D44" "0.6142"
Set-TextValue "E44" "  +1.28%  "
Set-TextValue "D45" "13.23"
Set-TextValue "E45" "  +1.20%  "
Set-TextValue "D46" "3.697"
Set-TextValue "E46" "  -0.11%  "
Set-TextValue "D47" "2.068"
Set-TextValue "E47" "  +1.89%  "
Set-TextValue "E48" "  +1.64%  "
Set-TextValue "D49" "124.78"
Set-TextValue "E49" "  +1.94%  "
Set-TextValue "D50" "1.187"
Set-TextValue "E50" "  -0.64%  "
Set-TextValue "D51" "78.55"
Set-TextValue "E51" "  +2.05%  "
